$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 91403.45
$ws.Range("I33").Value = 111569.445
$ws.Range("K33").Value = 111569.445
$ws.Range("M33").Value = -111340.445
$ws.Range("H88").Value = 6332.6665
$ws.Range("J88").Value = 6499.5
$ws.Range("L88").Value = 6499.5
$ws.Range("N88").Value = -7311.5
$ws.Range("H91").Value = 6332.6665
$ws.Range("J91").Value = 6499.5
$ws.Range("L91").Value = 6499.5
$ws.Range("N91").Value = -9307.5
$ws.Range("H141").Value = 5918.8
$ws.Range("I141").Value = 5918.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 17756.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -12576.4
$ws.Range("N141").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3878
$ws.Range("I2").Value = 2540.2
$ws.Range("K2").Value = 2540.2
$ws.Range("M2").Value = -2427.2
$ws.Range("H45").Value = 2011.8572
$ws.Range("J45").Value = 1991.25
$ws.Range("L45").Value = 1991.25
$ws.Range("N45").Value = -2745.25
$ws.Range("H54").Value = 44443
$ws.Range("J54").Value = 44443
$ws.Range("L54").Value = 44443
$ws.Range("N54").Value = -45981
$ws.Range("H116").Value = 3878
$ws.Range("I116").Value = 2540.2
$ws.Range("K116").Value = 2540.2
$ws.Range("M116").Value = -246.1999999999998

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3878
$ws.Range("I3").Value = 2540.2
$ws.Range("K3").Value = 2540.2
$ws.Range("M3").Value = -2426.2
$ws.Range("H86").Value = 1900
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -777
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1900
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3884
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 299.66666
$ws.Range("I94").Value = 223.4
$ws.Range("J94").Value = 395
$ws.Range("K94").Value = 223.4
$ws.Range("L94").Value = 395
$ws.Range("M94").Value = 227.6
$ws.Range("N94").Value = -1297
$ws.Range("H99").Value = 4742
$ws.Range("I99").Value = 4734
$ws.Range("K99").Value = 4734
$ws.Range("M99").Value = -3236

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2345.8
$ws.Range("I31").Value = 1463.3636
$ws.Range("J31").Value = 4772.5
$ws.Range("K31").Value = 1463.3636
$ws.Range("L31").Value = 4772.5
$ws.Range("M31").Value = -1168.3636
$ws.Range("N31").Value = -5362.5
$ws.Range("H34").Value = 2345.8
$ws.Range("I34").Value = 1463.3636
$ws.Range("J34").Value = 4772.5
$ws.Range("K34").Value = 1463.3636
$ws.Range("L34").Value = 4772.5
$ws.Range("M34").Value = -1261.3636
$ws.Range("N34").Value = -5176.5
$ws.Range("H50").Value = 32700
$ws.Range("I50").Value = 32700
$ws.Range("K50").Value = 32700
$ws.Range("M50").Value = -32075
$ws.Range("H55").Value = 50000
$ws.Range("I55").Value = 50000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 50000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -49685
$ws.Range("N55").ClearContents()
$ws.Range("H99").Value = 3062.0952
$ws.Range("I99").Value = 2673
$ws.Range("J99").Value = 4307.2
$ws.Range("K99").Value = 2673
$ws.Range("L99").Value = 4307.2
$ws.Range("M99").Value = -1175
$ws.Range("N99").Value = -7303.2
$ws.Range("H115").Value = 32700
$ws.Range("J115").Value = 32700
$ws.Range("L115").Value = 32700
$ws.Range("N115").Value = -35050
$ws.Range("H126").Value = 3062.0952
$ws.Range("I126").Value = 2673
$ws.Range("J126").Value = 4307.2
$ws.Range("K126").Value = 8019
$ws.Range("L126").Value = 12921.6
$ws.Range("M126").Value = -5549
$ws.Range("N126").Value = -17861.6

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1867.3334
$ws.Range("I55").Value = 3002
$ws.Range("K55").Value = 9006
$ws.Range("M55").Value = -8829
$ws.Range("H129").Value = 533.2857
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 20016.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 20016.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 20016.5
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -20612.5
$ws.Range("H52").Value = 15537.25
$ws.Range("J52").Value = 16033
$ws.Range("L52").Value = 16033
$ws.Range("N52").Value = -16551
$ws.Range("H80").Value = 3166
$ws.Range("H83").Value = 3166

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 825.2222
$ws.Range("I22").Value = 388.33334
$ws.Range("J22").Value = 1699
$ws.Range("K22").Value = 388.33334
$ws.Range("L22").Value = 1699
$ws.Range("M22").Value = -93.33334000000002
$ws.Range("N22").Value = -2289
$ws.Range("H27").Value = 825.2222
$ws.Range("I27").Value = 388.33334
$ws.Range("J27").Value = 1699
$ws.Range("K27").Value = 388.33334
$ws.Range("L27").Value = 1699
$ws.Range("M27").Value = -281.33334
$ws.Range("N27").Value = -1913
$ws.Range("H43").Value = 12419
$ws.Range("I43").Value = 9998
$ws.Range("J43").Value = 12605.23
$ws.Range("K43").Value = 9998
$ws.Range("L43").Value = 12605.23
$ws.Range("M43").Value = -9805
$ws.Range("N43").Value = -12991.23
$ws.Range("H46").Value = 6395.6
$ws.Range("J46").Value = 5744.5
$ws.Range("L46").Value = 5744.5
$ws.Range("N46").Value = -6120.5
$ws.Range("H82").Value = 5150.125
$ws.Range("I82").Value = 4116.5
$ws.Range("J82").Value = 8251
$ws.Range("K82").Value = 4116.5
$ws.Range("L82").Value = 8251
$ws.Range("M82").Value = -3755.5
$ws.Range("N82").Value = -8973
$ws.Range("H85").Value = 5150.125
$ws.Range("I85").Value = 4116.5
$ws.Range("J85").Value = 8251
$ws.Range("K85").Value = 4116.5
$ws.Range("L85").Value = 8251
$ws.Range("M85").Value = -2868.5
$ws.Range("N85").Value = -10747
$ws.Range("H92").Value = 28389
$ws.Range("J92").Value = 28389
$ws.Range("L92").Value = 28389
$ws.Range("N92").Value = -33381
$ws.Range("H100").Value = 3665.4546
$ws.Range("J100").Value = 4235.4
$ws.Range("L100").Value = 4235.4
$ws.Range("N100").Value = -5317.4

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 14602.286
$ws.Range("I37").Value = 19999.5
$ws.Range("J37").Value = 12443.4
$ws.Range("K37").Value = 19999.5
$ws.Range("L37").Value = 12443.4
$ws.Range("M37").Value = -19796.5
$ws.Range("N37").Value = -12849.4
$ws.Range("H62").Value = 3319.8
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 3319.8
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H81").Value = 9358.385
$ws.Range("I81").Value = 9358.385
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 18716.77
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -17655.77
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 9358.385
$ws.Range("I84").Value = 9358.385
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 93583.85000000001
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -88279.85000000001
$ws.Range("N84").ClearContents()
$ws.Range("H96").Value = 2708.7
$ws.Range("I96").Value = 4597.3335
$ws.Range("K96").Value = 4597.3335
$ws.Range("M96").Value = -3224.3335
$ws.Range("H122").Value = 7614.4165
$ws.Range("I122").Value = 7142.6
$ws.Range("K122").Value = 21427.8
$ws.Range("M122").Value = -18977.8
$ws.Range("H126").Value = 2185.9167
$ws.Range("I126").Value = 1520.2222
$ws.Range("J126").Value = 4183
$ws.Range("K126").Value = 4560.6666
$ws.Range("L126").Value = 12549
$ws.Range("M126").Value = -2090.6666
$ws.Range("N126").Value = -17489
